$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "Business Exceptions" sheet: keep only the first block of rows
#    (A2:A22) and drop the two duplicated blocks that used to follow
#    it (rows 23-62 in the original file).
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Business Exceptions")

$s4 = "CHANGES - SOX...eml on date01/24/2020 00:00:00 is missing ITRCA Member that filed or reviewed final evidence"
$s5 = '{row.item(3).ToString + " " + row.item(0).ToString.Substring(0,10) +  " does not contain Server Name " + row.item(7).ToString}'

# Row 17 used to hold the first "s5" value; it now holds one more "s4"
# value (matching what used to be row 37 of the second block).
$ws1.Range("A17").Value = $s4
# Row 21 used to be blank; row 22 used to start the next block with an
# "s1" value. Both now hold the "s5" value, extending that run to five
# rows (18-22) the same way the old second block did (rows 38-41 plus
# one extra row 22).
$ws1.Range("A21").Value = $s5
$ws1.Range("A22").Value = $s5

# Drop everything below the now-complete first block.
$ws1.Range("A23:A62").EntireRow.Delete()

# ------------------------------------------------------------------
# 2. "System Exceptions" sheet: remove the old sheet and replace it
#    with a brand-new, empty sheet of the same name (placed right
#    after "Sheet1", i.e. at the end of the tab strip).
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("System Exceptions")
$ws3.Delete()

$sheet1 = $wb.Worksheets.Item("Sheet1")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$newSheet.Name = "System Exceptions"

# A freshly inserted worksheet defaults to Excel's "blank workbook"
# page margins; restore the ones a brand-new sheet normally carries.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 54
$newSheet.PageSetup.BottomMargin = 36
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 54

# ------------------------------------------------------------------
# 3. Restore "Business Exceptions" as the active/selected tab.
# ------------------------------------------------------------------
$ws1.Activate()
